$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 6

$ws.Range("P10").Select()
